$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New version of the "parental relationship types" controlled vocabulary.
# Header row is renamed (lower-cased codes) and gains a third "note" column;
# a handful of rows are re-labelled and flagged as deprecated ("Non Valido")
# in the new note column.
# ---------------------------------------------------------------------------

# Row 1 - header
$ws.Cells.Item(1,1).Value = "codice_1_livello"
$ws.Cells.Item(1,2).Value = "label_1_livello_it"
$ws.Cells.Item(1,3).Value = "note"

# The header row used to be bold (style carried over from A2/B1); the new
# header row uses the regular (non-bold) font.
$ws.Cells.Item(1,1).Font.Bold = $false
$ws.Cells.Item(1,2).Font.Bold = $false
$ws.Cells.Item(1,3).Font.Bold = $false

# Row 2 - label rename, and it also loses its (already non-bold) explicit style
$ws.Cells.Item(2,2).Value = "Intestatario della Scheda"
$ws.Cells.Item(2,2).Font.Bold = $false

# Row 12 - casing fix: "Zio / Zia (Collaterale)" -> "Zio / Zia (collaterale)"
$ws.Cells.Item(12,2).Value = "Zio / Zia (collaterale)"

# Rows 15, 16, 20, 27 - newly deprecated entries (note = "Non Valido")
$ws.Cells.Item(15,3).Value = "Non Valido"
$ws.Cells.Item(16,3).Value = "Non Valido"
$ws.Cells.Item(20,3).Value = "Non Valido"
$ws.Cells.Item(27,3).Value = "Non Valido"

# Rows 21, 22 - casing fix: "(Affine)" -> "(affine)"
$ws.Cells.Item(21,2).Value = "Nipote (affine)"
$ws.Cells.Item(22,2).Value = "Zio / Zia (affine)"

# Rows 29, 30, 31 - also newly deprecated (note = "Non Valido")
$ws.Cells.Item(29,3).Value = "Non Valido"
$ws.Cells.Item(30,3).Value = "Non Valido"
$ws.Cells.Item(31,3).Value = "Non Valido"

# Widen the used range / selection to reflect the extra column and move the
# on-screen selection/scroll position to the new scratch cell below the table.
$excel.ActiveWindow.ScrollRow = 20
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C32").Select()
